$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (column C holds the "Förändrad" date)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 416 }

# Every cell in column C from row 2 through the last data row is updated
# from serial date 45180 to 45181 (i.e. incremented by one day).
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45181
